# Set line spacing to single (1x) for every paragraph in the document,
# matching <w:spacing w:line="240" w:lineRule="auto"/> in the OOXML.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.LineSpacingRule = 0   # wdLineSpaceSingle
}
